$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Charity/cause assessment ..." paragraph: split the run so the word
#    "neglectedness" is wrapped in spell-check proofErr markers, same as the
#    author's re-typed/re-checked text. Content stays the same, only the run
#    boundaries change.
# ---------------------------------------------------------------------------
$charityPara = $d.Paragraphs.Item(5)
$charityXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>
<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Charity/cause assessment &#8211; Scale, tractability, </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:strike/></w:rPr><w:t>neglectedness</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">, transparency, QALYs, </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$charityPara.Range.InsertXML($charityXml)

# ---------------------------------------------------------------------------
# 2. After "What can a student do to reduce suffering in the world?":
#      - drop the _GoBack bookmark that currently sits in that paragraph
#      - add a new "Guest blog" paragraph
#      - add two new blank paragraphs before "To DO:"
#    (InsertXML replaces the paragraph's own range *and* swallows the
#    paragraph mark that separated it from the next paragraph, so the
#    payload below must spell out both blank paragraphs explicitly.)
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$studentPara = $d.Paragraphs.Item(10)
$studentXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>What can a student do to reduce suffering in the world?</w:t></w:r></w:p>
<w:p><w:r><w:t>Guest blog</w:t></w:r></w:p>
<w:p/>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$studentPara.Range.InsertXML($studentXml)

# ---------------------------------------------------------------------------
# 3. Strike through the three "To DO" items (they're now crossed off).
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "Start EA assessment criteria blog post" -or $t -eq "Finish my pledge blog post" -or $t -eq "Finish EA page") {
        $p.Range.Font.StrikeThrough = 1
    }
}

# ---------------------------------------------------------------------------
# 4. Re-create the _GoBack bookmark at the end of "Finish my pledge blog
#    post" (that's where Word last left the cursor after these edits).
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq "Finish my pledge blog post") {
        $pledgePara = $p
    }
}
$pledgeXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:strike/></w:rPr></w:pPr>
<w:r><w:rPr><w:strike/></w:rPr><w:t>Finish my pledge blog post</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$pledgePara.Range.InsertXML($pledgeXml)
